{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// Change summary (from the diff):\n//   1. Remove bold (the <w:b> run property) from the four course\n//      section-header paragraphs:\n//        \"What is Happiness at Work and Why Does It Matter?\"\n//        \"How Can We Increase Our Own Happiness at Work?\"\n//        \"How to Be Happy With Others at Work\"\n//        \"How to Create a Happy Workplace\"\n//   2. Fill in the text of the empty bullet paragraph right before\n//      \"How to Be Happy With Others at Work\" with \"Progress Check 02\"\n//      (mirrors the \"Progress Check 01/03/04\" bullets elsewhere).\n//\n// Note: setting `paragraph.font.bold` (or `paragraph.getRange().font.bold`)\n// also stamps the paragraph-mark run properties (w:pPr/w:rPr), which the\n// target diff does not touch. Using `body.search(...)` to get a Range over\n// just the run text (no paragraph mark) and clearing bold there removes the\n// <w:b> element cleanly without adding a paragraph mark-run override.\n\nconst body = context.document.body;\n\nconst headings = [\n  \"What is Happiness at Work and Why Does It Matter?\",\n  \"How Can We Increase Our Own Happiness at Work?\",\n  \"How to Be Happy With Others at Work\",\n  \"How to Create a Happy Workplace\",\n];\n\nfor (const heading of headings) {\n  const results = body.search(heading, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const result of results.items) {\n    result.font.bold = false;\n  }\n}\n\nawait context.sync();\n\n// Locate the single empty-text bullet paragraph and give it the missing\n// \"Progress Check 02\" label.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text === \"\") {\n    paragraph.insertText(\"Progress Check 02\", Word.InsertLocation.replace);\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Change summary (from the diff):\n#   1. Remove bold (the <w:b> run property) from the four course\n#      section-header paragraphs:\n#        \"What is Happiness at Work and Why Does It Matter?\"\n#        \"How Can We Increase Our Own Happiness at Work?\"\n#        \"How to Be Happy With Others at Work\"\n#        \"How to Create a Happy Workplace\"\n#   2. Fill in the text of the empty bullet paragraph right before\n#      \"How to Be Happy With Others at Work\" with \"Progress Check 02\"\n#      (mirrors the \"Progress Check 01/03/04\" bullets elsewhere).\n#\n# Note: setting bold on a whole paragraph's Range (e.g. $p.Range.Font.Bold)\n# also stamps the paragraph-mark run properties (w:pPr/w:rPr), which the\n# target diff does not touch. Using Find.Execute() to locate just the run\n# text (no paragraph mark) and clearing bold on that found Range removes\n# the <w:b> element cleanly without adding a paragraph mark-run override.\n\n$d = $word.ActiveDocument\n\n$headings = @(\n    \"What is Happiness at Work and Why Does It Matter?\",\n    \"How Can We Increase Our Own Happiness at Work?\",\n    \"How to Be Happy With Others at Work\",\n    \"How to Create a Happy Workplace\"\n)\n\nforeach ($heading in $headings) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Text = $heading\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.Execute() | Out-Null\n    if ($range.Find.Found) {\n        $range.Font.Bold = 0\n    }\n}\n\n# Locate the single empty-text bullet paragraph and give it the missing\n# \"Progress Check 02\" label.\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd(\"`r\", \"`n\", \"`x07\")\n    if ($text -eq \"\") {\n        $p.Range.Text = \"Progress Check 02\"\n        break\n    }\n}\n"}
